$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before row 318. Excel shifts rows 318:429 down to
# 319:430 (preserving all their existing values/formatting), and the used
# range / dimension grows from A1:R429 to A1:R430 automatically.
$ws.Rows(318).Insert()

# Populate the newly-inserted row 318 with the new weekly record.
$ws.Range("A318").Value = 7
$ws.Range("B318").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C318").Value = "Ñuble"
$ws.Range("D318").Value = 44809
$ws.Range("E318").Value = 16
$ws.Range("F318").Value = 100114014
$ws.Range("G318").Value = "Betarraga"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 300
$ws.Range("K318").Value = 900
$ws.Range("L318").Value = 1000
$ws.Range("M318").Value = 950
$ws.Range("N318").Value = "`$/paquete 5 unidades"
$ws.Range("O318").Value = "Provincia de Diguillín"
$ws.Range("P318").Value = 190
$ws.Range("Q318").Value = 5
$ws.Range("R318").Value = "Hortaliza"
